$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the timestamp recorded for the existing last row (103) - the R script
# re-ran and produced a corrected intraday timestamp for that date.
$ws.Cells.Item(103, 1).Value = 45484.2916666667

# Append the new row (104) pulled in by the latest R script run.
$ws.Cells.Item(104, 1).Value = 45485.6176041667
$ws.Cells.Item(104, 2).Value = 3000
$ws.Cells.Item(104, 3).Value = 6.05999994277954
$ws.Cells.Item(104, 4).Value = 5.88000011444092
$ws.Cells.Item(104, 5).Value = 6
$ws.Cells.Item(104, 6).Value = 6.05999994277954

# Column G (adj_close) is stored as text in this sheet, so force the new value
# to be written as a text cell (matching the other rows) rather than a number -
# build it as a formula that returns the text, then freeze it down to a value.
$gCell = $ws.Cells.Item(104, 7)
$gCell.Formula = '="6.05999994277954"'
$gCell.Copy()
$gCell.PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Cells.Item(104, 8).Value = "PAL.MI"

# Give the new date cell (A104) the same date/time style used by the rest of
# column A (copy formats only, so the shared style is reused instead of a new
# one being minted).
$ws.Cells.Item(103, 1).Copy()
$ws.Cells.Item(104, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wb.Save()
